# Atualização 10/07 - Correção de variáveis, limpando imports inúteis, etc...
# Adds a new timesheet entry (row 4) to the "Ponto Eletrônico" sheet:
#   DATA=09/07/2023, ENTRADA=15:34:29, INTERVALO=15:34:30,
#   RETORNO INTERVALO=15:34:31, SAÍDA=15:34:33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as plain text (e.g. "09/07/2023"), matching
# the existing rows. Force text formatting first so Excel's autodetect
# doesn't silently turn the string into a real date serial number, then
# drop the temporary formatting so the cell keeps the sheet's default style.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "09/07/2023"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "15:34:29"
$ws.Range("C4").Value = "15:34:30"
$ws.Range("D4").Value = "15:34:31"
$ws.Range("E4").Value = "15:34:33"
